$wb = $excel.ActiveWorkbook

# Insert a new sheet; Worksheets.Add() drops it in before the active sheet,
# so it takes over index 1 / rId1 and the original sheet slides to index 2 / rId2
# (this matches the target: "modules" gets the new sheetId, "opportunities"
# keeps the original sheet's identity/sheetId).
$wsModules = $wb.Worksheets.Add()
$wsOpportunities = $wb.Worksheets.Item(2)

# The original sheet had three columns of data (Sale 1/2/3 across A1:C1);
# clear it out so we can rebuild it as a single column.
$wsOpportunities.Cells.Clear()

$wsModules.Name = "modules"
$wsOpportunities.Name = "opportunities"

# Fill in "opportunities" first so new shared strings land in the same order
# as the source workbook (opportunity before module name/CRM/Discuss/Calendar).
$wsOpportunities.Range("A1").Value = "opportunity"
$wsOpportunities.Range("A2").Value = "Sale 1"
$wsOpportunities.Range("A3").Value = "Sale 2"
$wsOpportunities.Range("A4").Value = "Sale 3"
$null = $wsOpportunities.Range("F32").Select()

# Fill in "modules".
$wsModules.Range("A1").Value = "module name"
$wsModules.Range("A2").Value = "CRM"
$wsModules.Range("A3").Value = "Discuss"
$wsModules.Range("A4").Value = "Calendar"

# Column A on "modules" is widened.
$wsModules.Columns.Item(1).ColumnWidth = 13.8

# Match the original sheet's page margins (values are in points: 72pt = 1in).
$wsModules.PageSetup.LeftMargin = 54
$wsModules.PageSetup.RightMargin = 54
$wsModules.PageSetup.TopMargin = 72
$wsModules.PageSetup.BottomMargin = 72
$wsModules.PageSetup.HeaderMargin = 36
$wsModules.PageSetup.FooterMargin = 36

$null = $wsModules.Range("D15").Select()
